# Fill in the already-present-but-empty test rows (9-13) with their input
# values and output formulas, then append two brand-new test case rows
# (14, 15) - "3 new added tests, as well as the given" per the commit
# message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 9: Green Package 2 months ----
$ws.Range("B9").Value = "Green"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0
$ws.Range("F9").Formula = "=C9*C2"
$ws.Range("G9").Formula = "=(IF(F9>75,F9-20,0))"

# ---- Row 10: Blue Package 2 months ----
$ws.Range("B10").Value = "Blue"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0
$ws.Range("F10").Formula = "=C10*C3"
$ws.Range("G10").Formula = "=C10*C3"

# ---- Row 11: Purple Package 2 months ----
$ws.Range("B11").Value = "Purple"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0
$ws.Range("F11").Formula = "=C11*C4"
$ws.Range("G11").Formula = "=C11*C4"

# ---- Row 12: Green Package 3 months with additional 5 GB data ----
$ws.Range("B12").Value = "Green"
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 5
$ws.Range("F12").Formula = "=C12*C2+(D12*D2)"
$ws.Range("G12").Formula = "=(IF(F12>75,F12-20,0))"

# ---- Row 13: Purple Package 4 months (new test case, row already existed blank) ----
$ws.Range("A13").Value = "Purple Package 4 months"
$ws.Range("B13").Value = "Purple"
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 0
$ws.Range("F13").Formula = "=C4*C13"
$ws.Range("G13").Formula = "=C4*C13"
$ws.Rows(13).RowHeight = 43.5

# ---- Row 14: Green Package 4 months with additional 3 GB data (new row) ----
$ws.Range("A9:D9").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("A14").Value = "Green Package 4 months with additional 3 GB data"
$ws.Range("B14").Value = "Green"
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 3
$ws.Range("F14").Formula = "=C14*C2+(D14*D2)"
$ws.Range("G14").Formula = "=(IF(F14>75,F14-20,0))"
$ws.Range("F14").HorizontalAlignment = -4131
$ws.Rows(14).RowHeight = 47.25

# ---- Row 15: Blue Package 3 months, with additional 2 GB data (new row) ----
$ws.Range("A9:D9").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)

$ws.Range("B15").Value = "Blue "
$ws.Range("A15").Value = "Blue Package 3 months, with additional 2 GB data"
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 2
$ws.Range("F15").Formula = "=C15*C3+(D15*D3)"
$ws.Range("G15").Formula = "=C15*C3+(D15*D3)"
$ws.Range("F15").HorizontalAlignment = -4131
$ws.Range("G15").HorizontalAlignment = -4131
$ws.Rows(15).RowHeight = 52.5

# ---- Header-row height tweak that came along with the new rows ----
$ws.Rows(8).RowHeight = 30.95

# ---- View state: leave selection on the last-edited cell ----
$ws.Range("G13").Select()
